{"js": "// Template placeholders get renamed:\n//   {surname} -> {surName}   (the lone \"n\" becomes \"N\")\n//   {jmbg}    -> {jmbgNum}   (\"Num\" gets appended right after \"jmbg\")\n// Scope every search to matchCase so we never touch the \"JMBG\" label text,\n// the already-correct \"{surName}\" occurrence further down in the document,\n// or any of the many other \"name\"/\"ame\" substrings used elsewhere.\n\nconst body = context.document.body;\n\n// --- Change 1: \"surname\" -> \"surName\" ---\nconst surnameMatches = body.search(\"surname\", { matchCase: true, matchWholeWord: false });\nsurnameMatches.load(\"items\");\nawait context.sync();\n\nif (surnameMatches.items.length > 0) {\n  const surnameRange = surnameMatches.items[0];\n  // \"surname\" contains exactly one \"n\"; search scoped to this range only\n  // finds that single character so neighbouring runs/paragraphs are untouched.\n  const nMatches = surnameRange.search(\"n\", { matchCase: true });\n  nMatches.load(\"items\");\n  await context.sync();\n\n  if (nMatches.items.length > 0) {\n    nMatches.items[0].insertText(\"N\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// --- Change 2: \"jmbg\" -> \"jmbgNum\" ---\nconst jmbgMatches = body.search(\"jmbg\", { matchCase: true, matchWholeWord: false });\njmbgMatches.load(\"items\");\nawait context.sync();\n\nif (jmbgMatches.items.length > 0) {\n  const jmbgRange = jmbgMatches.items[0];\n  const afterJmbg = jmbgRange.getRange(\"End\");\n  afterJmbg.insertText(\"Num\", \"Before\");\n  await context.sync();\n}\n", "ps1": "# Template placeholders get renamed:\n#   {surname} -> {surName}   (the lone \"n\" becomes \"N\")\n#   {jmbg}    -> {jmbgNum}   (\"Num\" gets appended right after \"jmbg\")\n# Every Find is case-sensitive and scoped to a duplicated range so we never\n# touch the \"JMBG\" label text, the already-correct \"{surName}\" occurrence\n# further down in the document, or any of the other \"name\" placeholders.\n\n$d = $word.ActiveDocument\n\nfunction Find-ScopedRange($scopeRange, $text) {\n  $r = $scopeRange.Duplicate\n  $r.Find.ClearFormatting()\n  $r.Find.Text = $text\n  $r.Find.MatchCase = $true\n  $r.Find.MatchWholeWord = $false\n  $r.Find.MatchWildcards = $false\n  $r.Find.Forward = $true\n  $r.Find.Wrap = 0\n  if ($r.Find.Execute()) {\n    return $r\n  }\n  return $null\n}\n\n# --- Change 1: \"surname\" -> \"surName\" ---\n$surnameRange = Find-ScopedRange $d.Content \"surname\"\nif ($surnameRange -ne $null) {\n  $nRange = Find-ScopedRange $surnameRange \"n\"\n  if ($nRange -ne $null) {\n    $nRange.Text = \"N\"\n  }\n}\n\n# --- Change 2: \"jmbg\" -> \"jmbgNum\" ---\n$jmbgRange = Find-ScopedRange $d.Content \"jmbg\"\nif ($jmbgRange -ne $null) {\n  $endPt = $d.Range($jmbgRange.End, $jmbgRange.End)\n  $endPt.InsertAfter(\"Num\")\n}\n"}
